$d = $word.ActiveDocument

function Insert-ParaXml($range, $innerXml) {
    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + $innerXml + '</w:document></pkg:xmlData></pkg:part></pkg:package>'
    $range.InsertXML($xml)
}

# ---------------------------------------------------------------------
# 1) Remove the paragraph-mark run formatting (rFonts hint=eastAsia) on
#    the "The data format " list paragraph's pPr.
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("The data format") | Out-Null
$para = $rng.Paragraphs.Item(1)
$pRange = $para.Range
Insert-ParaXml $pRange '<w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:ind w:firstLineChars="0"/></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve">The data format </w:t></w:r></w:p></w:body>'

# ---------------------------------------------------------------------
# 2) The empty paragraph right after the first table loses its pPr/rPr,
#    becoming a bare <w:p/>.
# ---------------------------------------------------------------------
$t1 = $d.Tables.Item(1)
$afterT1 = $d.Range($t1.Range.End, $t1.Range.End + 1)
Insert-ParaXml $afterT1 '<w:body><w:p/></w:body>'

# ---------------------------------------------------------------------
# 3) Change the trailing run " " -> ": First byte" on the
#    "The bit gramma of Action " paragraph, without touching the other
#    runs in that paragraph.
# ---------------------------------------------------------------------
$rng2 = $d.Content
$rng2.Find.Execute("gramma") | Out-Null
$rng2.Expand(4) | Out-Null  ## wdParagraph
$lastRunRange = $d.Range($rng2.End - 2, $rng2.End - 1)
Insert-ParaXml $lastRunRange '<w:body><w:p><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>: First byte</w:t></w:r></w:p></w:body>'

# ---------------------------------------------------------------------
# 4) The last 9-column table loses its first column ("9-16B" / blank),
#    and the remaining 8 columns get new widths.  Also the "Absolute M"
#    + "ove" runs merge into a single "Absolute Move" run.
# ---------------------------------------------------------------------
$t3 = $d.Tables.Item(3)
$t3.Columns.Item(1).Delete()

$widths = @(947, 947, 947, 949, 971, 950, 950, 949)
for ($c = 1; $c -le $t3.Columns.Count; $c++) {
    $t3.Columns.Item($c).Width = $widths[$c - 1] / 20.0
}

$moveCell = $t3.Cell(2, 5)
$moveParaRange = $moveCell.Range.Paragraphs.Item(1).Range
Insert-ParaXml $moveParaRange '<w:body><w:p><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>Absolute Move</w:t></w:r></w:p></w:body>'
